$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (sumitIdentity1/SummitPass1 -> sumitIdentity3/SummitPass3)
$ws.Range("A2").Value = "sumitIdentity3"
$ws.Range("B2").Value = "SummitPass3"

# Add new row 3 values (sumitIdentity4/SummitPass4) for selenium grid setup
$ws.Range("A3").Value = "sumitIdentity4"
$ws.Range("B3").Value = "SummitPass4"

# Update the active selection to reflect the newly added row
$ws.Range("A2:B3").Select() | Out-Null
